$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper technique notes:
#  - The engine silently merges two adjacent runs that end up with
#    identical run formatting once a part is re-serialized. Toggling a
#    (no-op) formatting property on/off over a range forces it to stay
#    in its own <w:r>, without altering its visible formatting.
#  - Appending text with Range.InsertAfter() merges best (keeps proper
#    <w:rPr>) when the target run is still a single contiguous run, so
#    new sentences are appended *before* any run-splitting is done.
#  - Document.Bookmarks.Add(name, range) only honours the given range
#    when Start<End; a collapsed (zero-length) range is mishandled and
#    the bookmark lands at document position 0. To plant a true
#    zero-width bookmark at an arbitrary spot, a 1-character spacer is
#    inserted, wrapped with the bookmark, then deleted again -- Word
#    keeps bookmarks pinned in place when their wrapped text shrinks to
#    nothing.
# ---------------------------------------------------------------------

# -----------------------------------------------------------------
# 1) "1 clock tick = 6.25x10^(-5)ms. " paragraph -> append new sentence
#    as its own run.
# -----------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("1 clock tick = 6.25x10^(-5)ms. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insStart = $target.End
$insText = " You can use this as a conversion to go back and forth between units."
$ip = $d.Range($insStart, $insStart)
$ip.InsertAfter($insText)
$newRange = $d.Range($insStart, $insStart + $insText.Length)
$newRange.Bold = 1
$newRange.Bold = 0

# -----------------------------------------------------------------
# 2) "...R for row, or a C for column of the first character." paragraph
#    -> append new sentence as its own run.
# -----------------------------------------------------------------
$target2 = $d.Content
$target2.Find.Execute("or a C for column of the first character.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insStart2 = $target2.End
$insText2 = " These pins were connected correspondingly to the pins outlined on the directions."
$ip2 = $d.Range($insStart2, $insStart2)
$ip2.InsertAfter($insText2)
$newRange2 = $d.Range($insStart2, $insStart2 + $insText2.Length)
$newRange2.Bold = 1
$newRange2.Bold = 0

# -----------------------------------------------------------------
# 3) Prescaler paragraph: append the new explanatory sentence, then
#    split the run after "...more easily", then move the "_GoBack"
#    bookmark from the final paragraph to the end of this paragraph.
# -----------------------------------------------------------------
$prescalerParaIndex = 9

# 3a) Append the extra sentence at the very end of the paragraph while
#     it is still one run, so it merges in cleanly with a proper rPr.
$para = $d.Paragraphs($prescalerParaIndex)
$appendPos = $para.Range.End - 1
$extraText = " We would be able to achieve lower frequencies since a scalar would let us lower our dividing number and circumvent the limitation we have on the size of the timer counters."
$ip3 = $d.Range($appendPos, $appendPos)
$ip3.InsertAfter($extraText)

# 3b) Split the run right after "...more easily" (before " and achieve
#     lower frequencies...") by nudging formatting across the tail.
$splitAnchor = $d.Content
$splitAnchor.Find.Execute("set clock ticks more easily", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $splitAnchor.End

$para2 = $d.Paragraphs($prescalerParaIndex)
$paraTextEnd = $para2.Range.End - 1

$afterRange = $d.Range($splitPos, $paraTextEnd)
$afterRange.Bold = 1
$afterRange.Bold = 0

# 3c) Relocate the "_GoBack" bookmark to the end of this paragraph
#     (adjacent bookmarkStart/bookmarkEnd, 0-width). Re-adding a
#     bookmark with the same name removes it from its previous spot
#     (the final paragraph), matching the diff.
$para3 = $d.Paragraphs($prescalerParaIndex)
$bmPos = $para3.Range.End - 1
$bmIp = $d.Range($bmPos, $bmPos)
$bmIp.InsertAfter("Z")
$bmWrap = $d.Range($bmPos, $bmPos + 1)
$d.Bookmarks.Add("_GoBack", $bmWrap)
$bmDel = $d.Range($bmPos, $bmPos + 1)
$bmDel.Text = ""

Write-Host "done"
